{"js": "// Anmerkung von Thomas hinzugef\u00fcgt\n//\n// 1) \"soviel\" -> \"so viel\" in the Maven intro paragraph.\n// 2) Insert \"DIE-interner \" before \"Tomcat-Server\" in the Mavenplugin\n//    paragraph (note added by Thomas about an internal Tomcat server).\n\nconst body = context.document.body;\n\n// --- Edit 1: \"soviel\" -> \"so viel\" -------------------------------------\nconst soviel = body.search(\"soviel\", { matchCase: true, matchWholeWord: true });\nsoviel.load(\"items\");\nawait context.sync();\n\nif (soviel.items.length > 0) {\n  soviel.items[0].insertText(\"so viel\", \"Replace\");\n  await context.sync();\n}\n\n// --- Edit 2: \"ein Tomcat-Server\" -> \"ein DIE-interner Tomcat-Server\" ----\nconst tomcat = body.search(\"ein Tomcat-Server\", { matchCase: true });\ntomcat.load(\"items\");\nawait context.sync();\n\nif (tomcat.items.length > 0) {\n  tomcat.items[0].insertText(\"ein DIE-interner Tomcat-Server\", \"Replace\");\n  await context.sync();\n}\n", "ps1": "# Anmerkung von Thomas hinzugef\u00fcgt\n#\n# 1) \"soviel\" -> \"so viel\" in the Maven intro paragraph.\n# 2) Insert \"DIE-interner \" before \"Tomcat-Server\" in the Mavenplugin\n#    paragraph (note added by Thomas about an internal Tomcat server).\n\n$d = $word.ActiveDocument\n\n# --- Edit 1: \"soviel\" -> \"so viel\" --------------------------------------\n$rng1 = $d.Content\n$rng1.Find.ClearFormatting()\n$rng1.Find.Execute(\"soviel\", $true, $true, $false, $false, $false, $true, 1, $false, \"so viel\", 2) | Out-Null\n\n# --- Edit 2: \"ein Tomcat-Server\" -> \"ein DIE-interner Tomcat-Server\" ----\n$rng2 = $d.Content\n$rng2.Find.ClearFormatting()\n$rng2.Find.Execute(\"ein Tomcat-Server\", $true, $false, $false, $false, $false, $true, 1, $false, \"ein DIE-interner Tomcat-Server\", 2) | Out-Null\n"}
